# QuickMapFacts_project_steps.docx edit script
$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Split the "To expedite..." paragraph's single run into three
#    runs, changing "unzipping the provided" -> "unzipping below
#    provided".
# ---------------------------------------------------------------
$target = $d.Paragraphs(16)
if ($target.Range.Text.TrimEnd([char]13) -ne "To expedite the process, you can skip to the deployment step by unzipping the provided eb-flask.zip file, which contains a pre-configured Flask application ready for deployment.") {
    throw "paragraph 16 text did not match expectation"
}

$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="347A9791" w14:textId="77777777" w:rsidR="006E1272" w:rsidRDefault="006E1272" w:rsidP="0081147C"><w:r><w:t xml:space="preserve">To expedite the process, you can skip to the deployment step by unzipping </w:t></w:r><w:r><w:t>below</w:t></w:r><w:r><w:t xml:space="preserve"> provided eb-flask.zip file, which contains a pre-configured Flask application ready for deployment.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.Range.InsertXML($xml1)

# ---------------------------------------------------------------
# 2) After that paragraph (and before the existing blank paragraph
#    that follows it) insert: blank para, URL para, blank para.
# ---------------------------------------------------------------
$target = $d.Paragraphs(16)
$insPoint = $d.Range($target.Range.End, $target.Range.End)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t>https://drive.google.com/file/d/1WNatH0vcjNGK0dvmxXrG-p_iTlBbNsdi/view?usp=drive_link</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint.InsertXML($xml2)

Write-Output "Steps 1-2 done"
